# Remove the last 4 node rows (ids "I","J","K","L") from the "nodes" sheet.
# Their A (id) and B (ALTURA) cells are reset to plain/empty cells, matching
# the style already used by the unused C/D/E columns on that sheet. Once
# those shared strings are no longer referenced anywhere, the rebuilt
# sharedStrings table drops them and every remaining string index shifts
# down accordingly (this is what re-points the "edges" sheet header cells
# from id1/id2/distancia's old indices to their new ones automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("nodes")

# Pick up the unused/empty style already applied to column C in these rows
# (border-only style, no fill/value) and stamp it onto A10:B13 so the cells
# keep that formatting instead of the id/number styles.
$ws.Range("C10").Copy()
$ws.Range("A10:B13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Now blank out the values so the cells become empty (style-only), which is
# what drops the id10..id13 / ALTURA values for those rows.
$ws.Range("A10:B13").ClearContents()
